$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of the last existing header cell (bold, centered,
# bordered) onto the three new header cells, then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team's season record (same W/L/T for every player row, 2-45)
for ($row = 2; $row -le 45; $row++) {
    $ws.Cells.Item($row, 30).Value = 78
    $ws.Cells.Item($row, 31).Value = 84
    $ws.Cells.Item($row, 32).Value = 0
}
